$d = $word.ActiveDocument

$replacements = @(
    @{old='280×4=1120'; new='170×9=1530'},
    @{old='965×8=7720'; new='827×7=5789'},
    @{old='887×5=4435'; new='387×7=2709'},
    @{old='208×7=1456'; new='106×3=318'},
    @{old='995×5=4975'; new='753×4=3012'},
    @{old='791×6=4746'; new='326×9=2934'},
    @{old='616×3=1848'; new='292×4=1168'},
    @{old='325×9=2925'; new='779×4=3116'},
    @{old='824×8=6592'; new='215×5=1075'},
    @{old='916×6=5496'; new='425×9=3825'},
    @{old='574×4=2296'; new='994×5=4970'},
    @{old='320×2=640';  new='950×4=3800'},
    @{old='957×3=2871'; new='424×2=848'},
    @{old='886×4=3544'; new='472×2=944'},
    @{old='596×7=4172'; new='168×3=504'},
    @{old='161×6=966';  new='676×7=4732'},
    @{old='477×7=3339'; new='932×6=5592'},
    @{old='200×4=800';  new='103×5=515'},
    @{old='533×4=2132'; new='206×9=1854'},
    @{old='772×8=6176'; new='353×6=2118'},
    @{old='266×6=1596'; new='261×9=2349'},
    @{old='214×8=1712'; new='702×2=1404'},
    @{old='453×2=906';  new='983×4=3932'},
    @{old='349×6=2094'; new='333×9=2997'},
    @{old='614×2=1228'; new='244×4=976'}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
